# Reproduce the "Download Add Users File" edit: the "email" column header
# becomes "email/user-id" (cell C1 on Sheet1). Re-touching the cell's style
# and moving the active selection mirrors the rest of the recorded diff
# (C1 picks up a distinct cell format, the sheet's selection ends on C4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in C1.
$ws.Range("C1").Value = "email/user-id"

# The edited cell is re-stamped with its own (still default-looking) style.
$ws.Range("C1").Style = "Normal"

# Leave the selection where the editor ended up after the edit.
$ws.Range("C4").Select() | Out-Null
